# Workbook/worksheet handles (per the harness, $excel.ActiveWorkbook is already the open workbook).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 20: fill in the hours worked + the new "Completed" note -----------
# (previously blank cells, already carrying style index 4 from the row above)
$ws.Range("B20").Value = 1
$ws.Range("C20").Value = "Update backend. No articleDetail class anymore."

# --- Row 27: extend the running SUM to include the newly-filled row 20 -----
$ws.Range("B27").Formula = "=SUM(B2:B20)"

# --- Window / view state (best effort - mirrors the author's scroll & ------
# selection position when they made the edit) -------------------------------
$win = $excel.ActiveWindow
$win.Left = 29580
$win.Top = -2500
$win.ScrollRow = 12
$win.ScrollColumn = 1

$ws.Range("C20").Select()
